# Fruta / hortaliza, semanal
#
# A new daily price record is inserted as row 25 (pushing the existing
# rows 25-96 down to 26-97). The sheet's used range grows from A1:T96
# to A1:T97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25, shifting rows 25..96 down to 26..97.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record's data.
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 44607
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100108
$ws.Cells.Item(25, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(25, 9).Value = 100108002
$ws.Cells.Item(25, 10).Value = "Mango"
$ws.Cells.Item(25, 11).Value = "Sin especificar"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 180
$ws.Cells.Item(25, 14).Value = 7000
$ws.Cells.Item(25, 15).Value = 7500
$ws.Cells.Item(25, 16).Value = 7222
$ws.Cells.Item(25, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(25, 18).Value = "Perú"
$ws.Cells.Item(25, 19).Value = 1806
$ws.Cells.Item(25, 20).Value = 4
